# Updated cryptos list (Price/Volume(1h) refresh + a WrappedEther/Chainlink
# row swap) as described by the commit's OOXML diff.
#
# Price-looking values (e.g. "222.99") are forced to Text via NumberFormat
# "@" before assignment and the cell Style is reset to "Normal" right after,
# so Excel's automatic number parsing doesn't turn them into floats (losing
# trailing zeros / precision) while also not leaving a stray custom style
# behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.607.96"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "1.793.79"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "222.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.21"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.33%  "
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0688"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0934"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("D12").Value = "2.050.19"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.785.27"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.49%  "
$ws.Range("D15").Value = "34.639.80"
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.629"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("E17").Value = "  +2.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").Value = "0.0₃0792"
$ws.Range("E20").Value = "  +6.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0516"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("D35").Value = "1.436.68"
$ws.Range("E35").Value = "  -4.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0190"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.76%  "
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "84.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.89%  "
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.909"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.53%  "
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0495"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.41%  "
$ws.Range("D47").Value = "1.950.01"
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.43%  "
